$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.080.95"
$ws.Range("E2").Value = "  -2.94%  "

$ws.Range("D3").Value = "2.253.60"
$ws.Range("E3").Value = "  -3.73%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.26%  "

$ws.Range("D5").Value = "'495.10"
$ws.Range("E5").Value = "  -1.68%  "

$ws.Range("D6").Value = "'127.12"
$ws.Range("E6").Value = "  -1.09%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'0.528"
$ws.Range("E8").Value = "  -1.30%  "

$ws.Range("D9").Value = "2.275.52"
$ws.Range("E9").Value = "  -3.11%  "

$ws.Range("E10").Value = "  -3.58%  "

$ws.Range("E11").Value = "  +0.27%  "

$ws.Range("D12").Value = "'0.322"
$ws.Range("E12").Value = "  +0.85%  "

$ws.Range("D13").Value = "'4.63"
$ws.Range("E13").Value = "  -3.34%  "

$ws.Range("D14").Value = "2.662.62"
$ws.Range("E14").Value = "  -3.40%  "

$ws.Range("D15").Value = "'21.61"
$ws.Range("E15").Value = "  +0.06%  "

$ws.Range("D16").Value = "54.036.02"
$ws.Range("E16").Value = "  -2.95%  "

$ws.Range("E17").Value = "  -1.94%  "

$ws.Range("D18").Value = "2.275.60"
$ws.Range("E18").Value = "  -2.08%  "

$ws.Range("D19").Value = "'9.93"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("E20").Value = "  +1.36%  "

$ws.Range("D21").Value = "'302.37"
$ws.Range("E21").Value = "  -2.33%  "

$ws.Range("D22").Value = "'6.33"
$ws.Range("E22").Value = "  +2.33%  "

$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("D24").Value = "'63.73"
$ws.Range("E24").Value = "  -2.39%  "

$ws.Range("E25").Value = "  +0.38%  "

$ws.Range("D26").Value = "'0.374"
$ws.Range("E26").Value = "  +1.05%  "

$ws.Range("D27").Value = "2.392.08"
$ws.Range("E27").Value = "  -2.39%  "

$ws.Range("D28").Value = "'0.147"
$ws.Range("E28").Value = "  +1.03%  "

$ws.Range("D29").Value = "'7.11"
$ws.Range("E29").Value = "  +0.52%  "

$ws.Range("D30").Value = "'164.89"
$ws.Range("E30").Value = "  -4.09%  "

$ws.Range("E31").Value = "  -2.25%  "

$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").Value = "'5.86"
$ws.Range("E32").Value = "  +1.45%  "

$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").Value = "0.0₃0677"
$ws.Range("E33").Value = "  -3.59%  "

$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("D35").Value = "'0.993"
$ws.Range("E35").Value = "  -0.35%  "

$ws.Range("E36").Value = "  +1.20%  "

$ws.Range("D37").Value = "'17.50"
$ws.Range("E37").Value = "  -0.57%  "

$ws.Range("D38").Value = "'1.18"
$ws.Range("E38").Value = "  +1.05%  "

$ws.Range("E39").Value = "  +6.59%  "

$ws.Range("D40").Value = "'3.63"
$ws.Range("E40").Value = "  -0.21%  "

$ws.Range("D41").Value = "'35.26"
$ws.Range("E41").Value = "  -2.30%  "

$ws.Range("E42").Value = "  +1.32%  "

$ws.Range("E43").Value = "  +1.27%  "

$ws.Range("E44").Value = "  -0.28%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'125.79"
$ws.Range("E45").Value = "  -0.36%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'4.82"
$ws.Range("E46").Value = "  +0.95%  "

$ws.Range("D47").Value = "'0.0889"
$ws.Range("E47").Value = "  -0.45%  "

$ws.Range("D48").Value = "'0.543"
$ws.Range("E48").Value = "  -1.95%  "

$ws.Range("D49").Value = "'237.65"
$ws.Range("E49").Value = "  +0.47%  "

$ws.Range("D50").Value = "'0.0479"
$ws.Range("E50").Value = "  +0.52%  "

$ws.Range("E51").Value = "  -0.42%  "
